$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------------------
# Shape 2: title textbox ("#32 CLIPBOARD IMAGE" / "for Neovim")
# ---------------------------------------------------------------------------
$titleShape = $s.Shapes.Item(2)
$titleTextRange = $titleShape.TextFrame.TextRange
$titleTextRange.Runs(1, 1).Text = "#31 NVIMS Update"
$titleTextRange.Runs(2, 1).Text = "  NeoVim Config Switcher "

# Resize width only; height stays the same as before the edit.
$titleShape.Width = 564.378346456693
$titleShape.Height = 122.42834855669291

# ---------------------------------------------------------------------------
# Shape 3: bullet list textbox ("Single solution:" / "Arch / Ubuntu" / ...)
# ---------------------------------------------------------------------------
$bodyShape = $s.Shapes.Item(3)
$bodyTextRange = $bodyShape.TextFrame.TextRange

# Drop the first two bullet paragraphs ("Single solution:" and
# "Arch / Ubuntu"), leaving the former 3rd/4th paragraphs
# ("Window 10 / 11" and "WSL2 + ") which become paragraphs 1 and 2.
$bodyTextRange.Paragraphs(1, 1).Delete()
$bodyTextRange.Paragraphs(1, 1).Delete()

$bodyTextRange2 = $bodyShape.TextFrame.TextRange
$bodyTextRange2.Paragraphs(1, 1).Text = "Move Repository"

$bodyTextRange3 = $bodyShape.TextFrame.TextRange
$bodyTextRange3.Paragraphs(2, 1).Text = "https://github.com/Traap/nvims"

# Flatten both remaining paragraphs to indent level 0 with no bullet.
$bodyTextRange4 = $bodyShape.TextFrame.TextRange
$bodyTextRange4.Paragraphs(1, 1).IndentLevel = 1
$bodyTextRange4.Paragraphs(1, 1).ParagraphFormat.Bullet.Type = 0
$bodyTextRange4.Paragraphs(2, 1).IndentLevel = 1
$bodyTextRange4.Paragraphs(2, 1).ParagraphFormat.Bullet.Type = 0

# Resize the shape to its new (smaller) extents.
$bodyShape.Width = 374.50795275590554
$bodyShape.Height = 69.62826971653543
